$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark. It originally spans from the
#    title paragraph down into the "La présentation de Florian ..."
#    paragraph, so just delete the bookmark itself (this removes both
#    the bookmarkStart and bookmarkEnd markers from the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Insert a new paragraph right after "Attention au raisonnement"
#    (and before the two trailing empty paragraphs) with the new
#    remark, and re-create the "_GoBack" bookmark (collapsed) at the
#    end of that paragraph's text.
$target = $null
$targetIndex = -1
$i = 1
foreach ($para in $d.Paragraphs) {
    $txt = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Attention au raisonnement") {
        $target = $para
        $targetIndex = $i
    }
    $i = $i + 1
}

$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Style = "Normal"

# Insert the text with a temporary trailing placeholder character so
# that, when we collapse a range right before it, the collapse point
# is a genuine "mid run" position rather than sitting exactly on the
# paragraph-end boundary.
$newPara.Range.Text = "Projet non terminé, en attente de changementX"

$bmPos = $newPara.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the temporary placeholder character now that the bookmark is
# anchored in place.
$phRange = $d.Range($newPara.Range.End - 2, $newPara.Range.End - 1)
$phRange.Delete()
